$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the primary_phone values for both firms with new phone numbers
# (kept as text so the leading "+" and digits are preserved verbatim).
$ws.Range("E2").Value = "'+420602456789"
$ws.Range("E3").Value = "'+420602654321"

# Remove the now-empty contacts/region/postal_code cells entirely
# (Clear removes the cell node altogether, matching the source edit).
$ws.Range("G2").Clear()
$ws.Range("J2").Clear()
$ws.Range("N2").Clear()
$ws.Range("G3").Clear()
$ws.Range("J3").Clear()
$ws.Range("N3").Clear()
